# Append five new simulation-log rows (16-20) to the "Data" sheet, mirroring
# the existing rows for the 2018.08.24 date / effective RS run.
#
# Columns: A Date | B Time | C Neuron Type | D Radius(nm) | E Thickness(um)
#          F Fdrive(kHz) | G Adrive(kPa) | H Tstim(ms) | I PRF(kHz)
#          J Duty factor | K Sim. Type | L # samples | M Comp. time(s)
#          N # spikes | O Latency(ms) | P Spike rate(sp/ms)
#
# Text-looking values (dates/times/labels) are written as a formula that
# evaluates to the literal string, then converted to a plain static value
# via Copy + PasteSpecial(xlPasteValues, -4163). Writing them straight to
# `.Value` would let Excel's automatic data-type recognition reinterpret
# strings such as "2018.08.24" or "19:57:16" as dates/times instead of text,
# which is not what the source log contains (plain shared strings).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($cell, [string]$text)
    $escaped = $text.Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

$rows = @(
    @{ Row = 16; Date = "2018.08.24"; Time = "19:57:16"; Neuron = "RS"; Radius = 32; Thickness = 0; Fdrive = 500; Adrive = 100; Tstim = 250; PRF = 0.1; Duty = 0.9500000000000001; SimType = "effective"; Samples = 6000; CompTime = 9.800000000000001; Spikes = 78; Latency = 37.1; SpikeRate = 0.4857362089357296 },
    @{ Row = 17; Date = "2018.08.24"; Time = "20:05:09"; Neuron = "RS"; Radius = 32; Thickness = 0; Fdrive = 500; Adrive = 100; Tstim = 250; PRF = 0.1; Duty = 0.91;                  SimType = "effective"; Samples = 6000; CompTime = 7.58;                  Spikes = 49; Latency = 38.55; SpikeRate = 0.4479769510540157 },
    @{ Row = 18; Date = "2018.08.24"; Time = "20:08:04"; Neuron = "RS"; Radius = 32; Thickness = 0; Fdrive = 500; Adrive = 100; Tstim = 250; PRF = 0.1; Duty = 0.9500000000000001; SimType = "effective"; Samples = 6000; CompTime = 9.49;                  Spikes = 78; Latency = 37.1;  SpikeRate = 0.4857362089357296 },
    @{ Row = 19; Date = "2018.08.24"; Time = "20:08:24"; Neuron = "RS"; Radius = 32; Thickness = 0; Fdrive = 500; Adrive = 100; Tstim = 250; PRF = 0.1; Duty = 0.92;                  SimType = "effective"; Samples = 6000; CompTime = 7.48;                  Spikes = 47; Latency = 38.15000000000001; SpikeRate = 0.4599823614175941 },
    @{ Row = 20; Date = "2018.08.24"; Time = "20:15:23"; Neuron = "RS"; Radius = 32; Thickness = 0; Fdrive = 500; Adrive = 100; Tstim = 250; PRF = 0.1; Duty = 0.96;                  SimType = "effective"; Samples = 6000; CompTime = 9.460000000000001; Spikes = 85; Latency = 36.8;  SpikeRate = 0.4878273340644669 }
)

foreach ($r in $rows) {
    $row = $r.Row

    Set-TextCell $ws.Cells.Item($row, 1) $r.Date
    Set-TextCell $ws.Cells.Item($row, 2) $r.Time
    Set-TextCell $ws.Cells.Item($row, 3) $r.Neuron

    $ws.Cells.Item($row, 4).Value = $r.Radius
    $ws.Cells.Item($row, 5).Value = $r.Thickness
    $ws.Cells.Item($row, 6).Value = $r.Fdrive
    $ws.Cells.Item($row, 7).Value = $r.Adrive
    $ws.Cells.Item($row, 8).Value = $r.Tstim
    $ws.Cells.Item($row, 9).Value = $r.PRF
    $ws.Cells.Item($row, 10).Value = $r.Duty

    Set-TextCell $ws.Cells.Item($row, 11) $r.SimType

    $ws.Cells.Item($row, 12).Value = $r.Samples
    $ws.Cells.Item($row, 13).Value = $r.CompTime
    $ws.Cells.Item($row, 14).Value = $r.Spikes
    $ws.Cells.Item($row, 15).Value = $r.Latency
    $ws.Cells.Item($row, 16).Value = $r.SpikeRate
}
